# EquivalencesPatrons.xlsx edit: ajout CommandManager, grille de correction
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Replace "ImageMenu" (E16) and "ImageManager" (E18) with "CommandManager"
$ws.Range("E16").Value = "CommandManager"
$ws.Range("E18").Value = "CommandManager"

# Update the active cell selection to J15 (as seen when correction grid was reviewed)
$ws.Range("J15").Select()
